# Update the date information on the cover/title slide.
#
# Slide 1, shape "TextBox 6" contains three paragraphs:
#   1) "2023 7th Symposium on AI-Embedded System-on-Chip"
#   2) "Kyungpook National University, Daehak-ro 80, No. 724, Daegu, Korea, Sept. 25, 2023"
#   3) "AI-Embedded Software-on-Chip (AI-SoC) Lab"
#
# This edit moves the event date from "Sept. 25, 2023" to "Oct. 28, 2022"
# and updates the leading year/ordinal from "2023 7th" to "2022 6th".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)          # "TextBox 6"
$tr = $sh.TextFrame.TextRange

# --- 1) "2023 7th ..." -> "2022 6th ..." ------------------------------
$full = $tr.Text
$old1 = "2023 7"
$idx1 = $full.IndexOf($old1)
if ($idx1 -ge 0) {
    $tr.Characters($idx1 + 1, $old1.Length).Text = "2022 6"
}

# --- 2) "... Korea, Sept. 25" -> "... Korea, Oct. " --------------------
$full = $tr.Text
$old2 = "Kyungpook National University, Daehak-ro 80, No. 724, Daegu, Korea, Sept. 25"
$idx2 = $full.IndexOf($old2)
if ($idx2 -ge 0) {
    $tr.Characters($idx2 + 1, $old2.Length).Text = "Kyungpook National University, Daehak-ro 80, No. 724, Daegu, Korea, Oct. "
}

# --- 3) ", 2023" -> "28, " then append a new "2022" run at the end ----
#        of the same paragraph (mirrors typing more text right after
#        the existing run, which PowerPoint keeps as its own run).
$full = $tr.Text
$old3 = ", 2023"
$idx3 = $full.IndexOf($old3)
if ($idx3 -ge 0) {
    $tr.Characters($idx3 + 1, $old3.Length).Text = "28, "
}

$para2 = $tr.Paragraphs(2, 1)
$para2.InsertAfter("2022")
